# Applies the commit's change:
#  1. Insert a new "Player Info" worksheet before the existing "ODI Batting"
#     sheet, with player metadata (ID, NAME, BATTING_HAND, BOWL_STYLE).
#  2. In "ODI Batting", rename the MATCH_CARD_LINK column to MATCH_CODE and
#     replace each full scorecard URL with just the trailing MatchCode value.

$wb = $excel.ActiveWorkbook
$originalFirst = $wb.Worksheets.Item(1)

# --- 1. New "Player Info" sheet, inserted before "ODI Batting" ---
$info = $wb.Worksheets.Add($originalFirst)
$info.Name = "Player Info"

$info.Cells.Item(1, 1).Value = "ID"
$info.Cells.Item(1, 2).Value = "NAME"
$info.Cells.Item(1, 3).Value = "BATTING_HAND"
$info.Cells.Item(1, 4).Value = "BOWL_STYLE"

$headerRow = $info.Range("A1:D1")
$headerRow.Font.Bold = $true
$headerRow.Borders.LineStyle = 1
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4160

$info.Cells.Item(2, 1).Value = "'5474"
$info.Cells.Item(2, 2).Value = "Rahmanullah Gurbaz"
$info.Cells.Item(2, 3).Value = "Right Handed"
$info.Cells.Item(2, 4).Value = "Does Not Bowl | Unknown"

# NOTE: worksheet handles captured before the Add() above become stale
# (they track the slot index, which the insertion shifted) so the
# "ODI Batting" sheet must be re-acquired by name afterwards.
$batting = $wb.Worksheets.Item("ODI Batting")

# --- 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE ---
$batting.Cells.Item(1, 4).Value = "MATCH_CODE"

$matchCodes = @("4444", "4446", "4448", "4525", "4528", "4530", "4537", "4538", "4539", "4582", "4585", "4588", "4671", "4674", "4675")

for ($i = 0; $i -lt $matchCodes.Count; $i++) {
    $row = $i + 2
    $batting.Cells.Item($row, 4).Value = "'" + $matchCodes[$i]
}

Write-Output "done"
